$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number and date range)
$ws.Range("A8").Value = "Volume 30   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/24/2023  Through  7/30/2023"

# Data table updates (rows 14-29)
$ws.Range("F14").Value = "'0"
$ws.Range("H14").Value = -100
$ws.Range("M14").Value = -38.095238095238
$ws.Range("N14").Value = -83.116883116883
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 9
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 38
$ws.Range("K15").Value = 31.034482758620
$ws.Range("L15").Value = 15.151515151515
$ws.Range("M15").Value = 5.555555555555
$ws.Range("N15").Value = -42.424242424242
$ws.Range("C16").Value = 19
$ws.Range("D16").Value = 15
$ws.Range("E16").Value = 26.666666666666
$ws.Range("F16").Value = 57
$ws.Range("G16").Value = 52
$ws.Range("H16").Value = 9.615384615384
$ws.Range("I16").Value = 393
$ws.Range("J16").Value = 460
$ws.Range("K16").Value = -14.565217391304
$ws.Range("L16").Value = 27.184466019417
$ws.Range("M16").Value = -5.301204819277
$ws.Range("N16").Value = -76.579261025029
$ws.Range("C17").Value = 24
$ws.Range("D17").Value = 25
$ws.Range("E17").Value = -4
$ws.Range("F17").Value = 105
$ws.Range("G17").Value = 95
$ws.Range("H17").Value = 10.526315789473
$ws.Range("I17").Value = 633
$ws.Range("J17").Value = 651
$ws.Range("K17").Value = -2.764976958525
$ws.Range("L17").Value = 30.785123966942
$ws.Range("M17").Value = 47.209302325581
$ws.Range("N17").Value = -27.241379310344
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = 37.5
$ws.Range("F18").Value = 37
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = 27.586206896551
$ws.Range("I18").Value = 230
$ws.Range("J18").Value = 231
$ws.Range("K18").Value = -0.432900432900
$ws.Range("L18").Value = 1.769911504424
$ws.Range("M18").Value = -15.129151291512
$ws.Range("N18").Value = -78.136882129277
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 81
$ws.Range("G19").Value = 114
$ws.Range("H19").Value = -28.947368421052
$ws.Range("I19").Value = 635
$ws.Range("J19").Value = 775
$ws.Range("K19").Value = -18.064516129032
$ws.Range("L19").Value = 8.177172061328
$ws.Range("M19").Value = 57.960199004975
$ws.Range("N19").Value = 22.115384615384
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 71.428571428571
$ws.Range("F20").Value = 50
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 315
$ws.Range("J20").Value = 286
$ws.Range("K20").Value = 10.139860139860
$ws.Range("L20").Value = 5.351170568561
$ws.Range("M20").Value = 77.966101694915
$ws.Range("N20").Value = -79.690522243713
$ws.Range("C21").Value = 92
$ws.Range("D21").Value = 79
$ws.Range("E21").Value = 16.455696202531
$ws.Range("F21").Value = 339
$ws.Range("G21").Value = 329
$ws.Range("H21").Value = 3.039513677811
$ws.Range("I21").Value = 2257
$ws.Range("J21").Value = 2443
$ws.Range("K21").Value = -7.613589848546
$ws.Range("L21").Value = 15.802975885069
$ws.Range("M21").Value = 28.824200913242
$ws.Range("N21").Value = -61.179910560715
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 26
$ws.Range("K22").Value = -48
$ws.Range("L22").Value = 30
$ws.Range("M22").Value = -35
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -44.444444444444
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = 33
$ws.Range("H23").Value = -39.393939393939
$ws.Range("I23").Value = 207
$ws.Range("J23").Value = 205
$ws.Range("K23").Value = 0.975609756097
$ws.Range("L23").Value = 24.698795180722
$ws.Range("M23").Value = 81.578947368421
$ws.Range("C24").Value = 44
$ws.Range("D24").Value = 73
$ws.Range("E24").Value = -39.726027397260
$ws.Range("F24").Value = 212
$ws.Range("G24").Value = 253
$ws.Range("H24").Value = -16.205533596837
$ws.Range("I24").Value = 1395
$ws.Range("J24").Value = 1691
$ws.Range("K24").Value = -17.504435245416
$ws.Range("L24").Value = 9.326018808777
$ws.Range("M24").Value = 62.020905923344
$ws.Range("C25").Value = 31
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 55
$ws.Range("F25").Value = 116
$ws.Range("G25").Value = 88
$ws.Range("H25").Value = 31.818181818181
$ws.Range("I25").Value = 741
$ws.Range("J25").Value = 768
$ws.Range("K25").Value = -3.515625
$ws.Range("L25").Value = 32.321428571428
$ws.Range("M25").Value = -31.893382352941
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 400
$ws.Range("F26").Value = 12
$ws.Range("G26").Value = 8
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 57
$ws.Range("J26").Value = 44
$ws.Range("K26").Value = 29.545454545454
$ws.Range("L26").Value = 0
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 400
$ws.Range("F27").Value = 12
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 71
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = 1.428571428571
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 14
$ws.Range("H28").Value = -78.571428571428
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = -30
$ws.Range("L28").Value = -4.545454545454
$ws.Range("M28").Value = -23.636363636363
$ws.Range("N28").Value = -81.081081081081
$ws.Range("C29").Value = "'0"
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 12
$ws.Range("H29").Value = -75
$ws.Range("J29").Value = 52
$ws.Range("K29").Value = -30.769230769230
$ws.Range("L29").Value = -7.692307692307
$ws.Range("M29").Value = -21.739130434782
$ws.Range("N29").Value = -81.909547738693
